# "added files from laptop"
# Adds three new replay-world entries (rows 2-4, column B) that had no
# name/title recorded yet, and leaves the selection on the next empty
# row the author was about to fill in (C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "unkown"
$ws.Range("B3").Value = "unkown"
$ws.Range("B4").Value = "unkown"

$ws.Range("C4").Select() | Out-Null
